$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text: add a colon after "Auszubildenden" in the big instructions
#    cell (A1), right before the #idName placeholder.
# ---------------------------------------------------------------------------
$headerText = $ws.Range("A1").Value2
$headerText = $headerText -replace "Auszubildenden #idName", "Auszubildenden: #idName"
$ws.Range("A1").Value2 = $headerText

# ---------------------------------------------------------------------------
# 2) "Lernfeld Nummer" column (D3) should use an integer number format
#    instead of a 2-decimal one -> mints a new style (numFmtId 1 = "0").
# ---------------------------------------------------------------------------
$ws.Range("D3").NumberFormat = "0"

# ---------------------------------------------------------------------------
# 3) Printing: fit the sheet to one page wide, scaled to 93%, unlimited tall.
# ---------------------------------------------------------------------------
$ws.PageSetup.Zoom = 93
$ws.PageSetup.FitToPagesWide = 1

# ---------------------------------------------------------------------------
# 4) Drop the workbook-level Print_Area defined name (now handled purely via
#    PageSetup fit-to-page above).
# ---------------------------------------------------------------------------
$ws.PageSetup.PrintArea = ""

# ---------------------------------------------------------------------------
# 5) Move the window a bit to the right and change the active selection.
# ---------------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = 1100

[void]$ws.Range("G1").Select()
